# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 24 and 25) at the top of the
# data block for "Terminal Hortofrutícola Agro Chillán - Membrillo",
# pushing the previously existing rows 24-43 down to rows 26-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 24 (Excel shifts rows 24:43 -> 26:45
# and copies formatting from the row immediately above, which already carries
# the date style used throughout column D).
$ws.Rows("24:25").Insert()

# ---- New row 24 : Especial ----
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 45096
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104003
$ws.Range("J24").Value = "Membrillo"
$ws.Range("K24").Value = "Champion"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 11000
$ws.Range("O24").Value = 11000
$ws.Range("P24").Value = 11000
$ws.Range("Q24").Value = "$/caja 18 kilos empedrada"
$ws.Range("R24").Value = "Región del Maule"
$ws.Range("S24").Value = 611
$ws.Range("T24").Value = 18

# ---- New row 25 : Primera ----
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 45096
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100104
$ws.Range("H25").Value = "Frutos de pepita"
$ws.Range("I25").Value = 100104003
$ws.Range("J25").Value = "Membrillo"
$ws.Range("K25").Value = "Champion"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 9500
$ws.Range("Q25").Value = "$/caja 18 kilos empedrada"
$ws.Range("R25").Value = "Región del Maule"
$ws.Range("S25").Value = 528
$ws.Range("T25").Value = 18
